$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column F
$ws.Range("F1").Value = "dmax_95CI"

# Add new dmax_95CI estimate + CI values for each data row (same value repeated)
$ws.Range("F2").Value = "0.4188 (0.2205 - 0.6171)"
$ws.Range("F3").Value = "0.4188 (0.2205 - 0.6171)"
$ws.Range("F4").Value = "0.4188 (0.2205 - 0.6171)"
